$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.194.23"
$ws.Range("E2").Value = "  +4.25%  "

$ws.Range("D3").Value = "2.490.58"
$ws.Range("E3").Value = "  +2.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.61"
$ws.Range("E5").Value = "  +1.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.92"
$ws.Range("E6").Value = "  +5.26%  "

$ws.Range("E7").Value = "  +2.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  +3.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.15"
$ws.Range("E10").Value = "  +7.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +1.87%  "

$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.40"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.19"

$ws.Range("D15").Value = "2.884.33"
$ws.Range("E15").Value = "  +2.85%  "

$ws.Range("D16").Value = "2.509.08"
$ws.Range("E16").Value = "  +3.88%  "

$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").Value = "47.126.32"
$ws.Range("E18").Value = "  +4.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.89"
$ws.Range("E19").Value = "  +5.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.69"
$ws.Range("E20").Value = "  +5.67%  "

$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.55"
$ws.Range("E22").Value = "  +2.29%  "

$ws.Range("E23").Value = "  +6.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.63"
$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("E25").Value = "  +4.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.14"
$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +4.02%  "

$ws.Range("E29").Value = "  -3.03%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  +11.06%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.91"
$ws.Range("E31").Value = "  +6.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.36"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.48"
$ws.Range("E33").Value = "  +5.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.58"
$ws.Range("E34").Value = "  -2.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0789"
$ws.Range("E35").Value = "  +3.45%  "

$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.98"
$ws.Range("E37").Value = "  +5.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.66"
$ws.Range("E38").Value = "  +5.15%  "

$ws.Range("E39").Value = "  +3.82%  "

$ws.Range("E40").Value = "  +1.97%  "

$ws.Range("E41").Value = "  +2.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.78"
$ws.Range("E42").Value = "  -3.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.24"
$ws.Range("E43").Value = "  +2.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0298"
$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("D45").Value = "1.964.02"
$ws.Range("E45").Value = "  +1.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.01"
$ws.Range("E46").Value = "  +1.94%  "

$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.79"
$ws.Range("E48").Value = "  +0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.05"
$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.28"
$ws.Range("E50").Value = "  +9.34%  "

$ws.Range("E51").Value = "  +3.70%  "
